$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("H 72") entirely; this shifts all following rows up by one,
# matching the diff (dimension A1:F63 -> A1:F62).
$ws.Rows.Item(2).Delete()
